# Automatische test-sync: 2025-08-05 19:35:50
# Adds the new "Testmail #7" log entry to the "Logs" sheet (row 48),
# extends the conditional formatting ranges to include it, and updates
# the "Dashboard" sheet category ordering/counts to reflect the new entry.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# 1. Logs sheet: append the new row 48 with the Testmail #7 data
# ---------------------------------------------------------------
$logs = $wb.Worksheets.Item("Logs")

$logs.Cells.Item(48, 1).Value = "Is dit artikel momenteel beschikbaar?"
$logs.Cells.Item(48, 2).Value = "mailmind.test@zohomail.eu"
$logs.Cells.Item(48, 3).Value = "Testmail #7: Is dit artikel momenteel beschikbaar?"
$logs.Cells.Item(48, 4).Value = "Overig"
$logs.Cells.Item(48, 5).Value = @"
Beste afzender,
Hartelijk dank voor uw interesse in ons artikel. Om u nauwkeurig te kunnen informeren over de beschikbaarheid, zou u ons meer details kunnen verstrekken, zoals de naam of het artikelnummer van het gewenste product? Op basis van deze informatie kan ik voor u nagaan of het artikel momenteel op voorraad is. 
Met vriendelijke groet,
[Naam]
E-mailassistent
"@
$logs.Cells.Item(48, 6).Value = "2025-08-05 19:35:34"
$logs.Cells.Item(48, 7).Value = "Ja"
$logs.Cells.Item(48, 8).Value = "Nee"
$logs.Cells.Item(48, 9).Value = "Ja"
$logs.Cells.Item(48, 10).Value = "Nee"

# ---------------------------------------------------------------
# 2. Logs sheet: extend conditional formatting ranges from row 47
#    to row 48 (one rule per block is enough, as the rules in a
#    block share the same sqref).
# ---------------------------------------------------------------
$logs.Range("D2:D47").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("D2:D48"))
$logs.Range("G2:G47").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("G2:G48"))
$logs.Range("H2:H47").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("H2:H48"))
$logs.Range("I2:I47").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("I2:I48"))
$logs.Range("J2:J47").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("J2:J48"))

# ---------------------------------------------------------------
# 3. Dashboard sheet: the new log entry is categorised "Overig",
#    which now ties "Inkoop / Bestellingen" at 6 occurrences.  The
#    category list is re-sorted (descending count, alphabetical on
#    ties), so rows 3-5 swap around:
#      row3: Inkoop / Bestellingen (6) -> Overig (6)
#      row4: Klantenservice / Contact (5) -> Inkoop / Bestellingen (6)
#      row5: Overig (5) -> Klantenservice / Contact (5)
# ---------------------------------------------------------------
$dash = $wb.Worksheets.Item("Dashboard")

$dash.Cells.Item(3, 1).Value = "Overig"
$dash.Cells.Item(4, 1).Value = "Inkoop / Bestellingen"
$dash.Cells.Item(4, 2).Value = 6
$dash.Cells.Item(5, 1).Value = "Klantenservice / Contact"
